# Updates cryptos list values (Price / Volume(1h) columns, plus a row swap
# between ShibaInu and WrappedEther) to match the latest scraped data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some "Price" strings (e.g. "1.059", "0.9712", "1.0000") look like
# plain numbers to Excel's auto-detection and would otherwise be silently
# converted to numeric values when assigned via .Value. Forcing the cell to
# text format, assigning, then restoring the "Normal" style keeps the value
# as literal text (matching the source data) without leaving a stray
# explicit number format on the cell.
function Set-TextValue {
    param($Coord, $Val)
    $cell = $ws.Range($Coord)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '21.757.21'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '1.539.96'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.07%  '
Set-TextValue "D6" '290.19'
$ws.Range("E6").Value = '  +0.67%  '
Set-TextValue "D7" '0.3885'
$ws.Range("E7").Value = '  +2.50%  '
$ws.Range("E8").Value = '  -2.50%  '
Set-TextValue "D9" '43.16'
$ws.Range("E9").Value = '  -0.17%  '
Set-TextValue "D10" '0.07206'
$ws.Range("E10").Value = '  -1.86%  '
Set-TextValue "D11" '1.059'
$ws.Range("E11").Value = '  -6.69%  '
Set-TextValue "D12" '1.001'
$ws.Range("E12").Value = '  -0.06%  '
Set-TextValue "D13" '5.638'
$ws.Range("E13").Value = '  -2.66%  '
Set-TextValue "D14" '18.60'
$ws.Range("E14").Value = '  -6.37%  '
Set-TextValue "D15" '6.606'
$ws.Range("E15").Value = '  -3.89%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.543.89'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D17" '0.00001113'
$ws.Range("E17").Value = '  +2.19%  '
Set-TextValue "D18" '0.06586'
$ws.Range("E18").Value = '  -1.01%  '
Set-TextValue "D19" '83.39'
$ws.Range("E19").Value = '  -2.53%  '
Set-TextValue "D20" '1.0000'
$ws.Range("E20").Value = '  -0.04%  '
Set-TextValue "D21" '6.146'
$ws.Range("E21").Value = '  -4.67%  '
Set-TextValue "D22" '15.40'
$ws.Range("E22").Value = '  -4.19%  '
Set-TextValue "D23" '10.98'
$ws.Range("E23").Value = '  -5.80%  '
Set-TextValue "D24" '2.388'
$ws.Range("E24").Value = '  +5.64%  '
$ws.Range("D25").Value = '21.764.08'
$ws.Range("E25").Value = '  -1.63%  '
Set-TextValue "D26" '2.385'
$ws.Range("E26").Value = '  -5.94%  '
Set-TextValue "D27" '146.68'
$ws.Range("E27").Value = '  -2.20%  '
Set-TextValue "D28" '18.38'
$ws.Range("E28").Value = '  -3.71%  '
Set-TextValue "D29" '4.848'
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("D30").Value = '1.718.42'
$ws.Range("E30").Value = '  -1.12%  '
Set-TextValue "D31" '117.48'
$ws.Range("E31").Value = '  -2.80%  '
Set-TextValue "D32" '0.9712'
$ws.Range("E32").Value = '  -12.80%  '
Set-TextValue "D33" '5.905'
$ws.Range("E33").Value = '  -1.46%  '
Set-TextValue "D34" '0.08212'
$ws.Range("E34").Value = '  +0.76%  '
Set-TextValue "D35" '8.913'
$ws.Range("E35").Value = '  -4.35%  '
Set-TextValue "D36" '5.154'
Set-TextValue "D37" '0.06081'
$ws.Range("E37").Value = '  -1.59%  '
Set-TextValue "D38" '1.488'
$ws.Range("E38").Value = '  -17.40%  '
Set-TextValue "D39" '0.02206'
$ws.Range("E39").Value = '  -3.59%  '
Set-TextValue "D40" '0.2039'
$ws.Range("E40").Value = '  -4.55%  '
Set-TextValue "D41" '1.191'
$ws.Range("E41").Value = '  -2.83%  '
$ws.Range("E42").Value = '  -0.09%  '
Set-TextValue "D43" '10.69'
$ws.Range("E43").Value = '  -2.71%  '
Set-TextValue "D44" '0.5757'
$ws.Range("E44").Value = '  -3.50%  '
Set-TextValue "D45" '3.747'
$ws.Range("E45").Value = '  +0.24%  '
Set-TextValue "D46" '12.98'
$ws.Range("E46").Value = '  -5.31%  '
Set-TextValue "D47" '0.5527'
$ws.Range("E47").Value = '  -4.10%  '
Set-TextValue "D48" '118.17'
$ws.Range("E48").Value = '  -1.65%  '
Set-TextValue "D49" '1.871'
$ws.Range("E49").Value = '  -4.80%  '
Set-TextValue "D50" '1.145'
$ws.Range("E50").Value = '  -1.93%  '
Set-TextValue "D51" '0.06729'
$ws.Range("E51").Value = '  -3.44%  '
